# Add a new "数组" (Array) worksheet at the end of the workbook containing
# LeetCode 26 "Remove Duplicates from Sorted Array" (solved with a hashmap),
# matching the layout/style of the existing sheets.

$wb = $excel.ActiveWorkbook

# Template sheet to copy header/data-row formatting (fonts, styles, etc.) from.
$template = $wb.Worksheets.Item(4)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Add the new worksheet after the last existing sheet.
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "数组"

# Copy the header row (row 1) and the first data row (row 2) formatting from
# the "数学" sheet so fonts / alignment / styles match the rest of the workbook.
$template.Range("A1:G1").Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)
$template.Range("A2:G2").Copy()
$ws.Range("A2:G2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths matching the committed sheet.
$ws.Columns.Item(1).ColumnWidth = 8.428571428571429
$ws.Columns.Item(2).ColumnWidth = 12.428571428571429
$ws.Columns.Item(3).ColumnWidth = 60.42857142857143
$ws.Columns.Item(4).ColumnWidth = 75.42857142857143
$ws.Columns.Item(5).ColumnWidth = 18.142857142857142
$ws.Columns.Item(6).ColumnWidth = 14.428571428571429
$ws.Columns.Item(7).ColumnWidth = 14.0

# Header row content ("No.", "leetcode", "题目", "解题方法", "解题关键词", "时间复杂度", "空间复杂度")
$ws.Cells.Item(1,1).Value = "No."
$ws.Cells.Item(1,2).Value = "leetcode"
$ws.Cells.Item(1,3).Value = "题目"
$ws.Cells.Item(1,4).Value = "解题方法"
$ws.Cells.Item(1,5).Value = "解题关键词"
$ws.Cells.Item(1,6).Value = "时间复杂度"
$ws.Cells.Item(1,7).Value = "空间复杂度"
$ws.Rows.Item(1).RowHeight = 22

$problemText = @'

给定一个排序数组，你需要在 原地 删除重复出现的元素，使得每个元素只出现一次，返回移除后数组的新长度。 
 不要使用额外的数组空间，你必须在 原地 修改输入数组 并在使用 O(1) 额外空间的条件下完成。
 示例 1:
 给定数组 nums = [1,1,2],
函数应该返回新的长度 2, 并且原数组 nums 的前两个元素被修改为 1, 2。
你不需要考虑数组中超出新长度后面的元素。
 示例 2:
 给定 nums = [0,0,1,1,1,2,2,3,3,4],
函数应该返回新的长度 5, 并且原数组 nums 的前五个元素被修改为 0, 1, 2, 3, 4。
你不需要考虑数组中超出新长度后面的元素。
 说明: 
 为什么返回数值是整数，但输出的答案是数组呢?
 请注意，输入数组是以「引用」方式传递的，这意味着在函数里修改输入数组对于调用者是可见的。 
 你可以想象内部操作如下:
 // nums 是以“引用”方式传递的。也就是说，不对实参做任何拷贝
int len = removeDuplicates(nums);
// 在函数里修改输入数组对于调用者是可见的。
// 根据你的函数返回的长度, 它会打印出数组中该长度范围内的所有元素。
for (int i = 0; i < len; i++) {
    print(nums[i]);
}
'@

# Data row: #1, leetcode 26, problem statement; solution/keyword/complexity left blank.
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 26
$ws.Cells.Item(2,3).Value = $problemText
$ws.Rows.Item(2).RowHeight = 409.6

# Selection/view state on the new sheet.
$ws.Range("D2").Select()

# Update the "数学" sheet's selection (no longer the scrolled-to C23 cell).
$math = $wb.Worksheets.Item(4)
$math.Range("A1:G2").Select()

# Make the new "数组" sheet the active tab.
$ws.Activate()
